$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LJ Speech")

# Row 2
$ws.Range("B2").Value = "<kero>"
$ws.Range("C2").Value = 24

# Row 3
$ws.Range("B3").Value = "<ant>"
$ws.Range("C3").Value = 30

# Row 4
$ws.Range("C4").Value = 30

# Row 5
$ws.Range("C5").Value = 41

# Row 6
$ws.Range("B6").Value = "<tab>"
$ws.Range("C6").Value = 37

# Row 7
$ws.Range("C7").Value = 36

# Row 9
$ws.Range("B9").Value = "<nine>"
$ws.Range("C9").Value = 33

# Row 10
$ws.Range("C10").Value = 32

# Row 11
$ws.Range("C11").Value = 32

# Row 12
$ws.Range("C12").Value = 32

# Row 13
$ws.Range("B13").Value = "<nike>"
$ws.Range("C13").Value = 33

# Row 15
$ws.Range("C15").Value = 12
